$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.169.29"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.144.80"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.49%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.00"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.31"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.72%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.134.94"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.163"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.91"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +5.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.458"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.54"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.664.86"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.27"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.949.59"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.138.34"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "468.44"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.42"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.32%  "
$ws.Range("E22").Value = "  +1.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.59"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.92%  "
$ws.Range("E24").Value = "  +12.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.14"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.05"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("E28").Value = "  +11.23%  "
$ws.Range("E29").Value = "  +1.60%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.28"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +6.24%  "
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("E33").Value = "  +3.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.70"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +4.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0857"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.73%  "
$ws.Range("E36").Value = "  +3.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.16"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.93%  "
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.28"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "461.47"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +6.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.27"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.91%  "
$ws.Range("E42").Value = "  +7.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.291"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +8.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0374"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.890.93"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.85"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +11.18%  "
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.53"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +8.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.111"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.23"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.83%  "
